# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# The "K" column (column G, header "K") is recomputed here with the newly
# calculated values for each data row (rows 2-46 on Sheet1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 1
    3  = 0
    4  = 2
    5  = 1
    6  = 1
    7  = 0
    8  = 1
    9  = 1
    10 = 0
    11 = 2
    12 = 1
    13 = 0
    14 = 1
    15 = 0
    16 = 1
    17 = 1
    18 = 2
    19 = 1
    20 = 3
    21 = 1
    22 = 0
    23 = 0
    24 = 0
    25 = 1
    26 = 1
    27 = 0
    28 = 0
    29 = 0
    30 = 2
    31 = 1
    32 = 1
    33 = 2
    34 = 1
    35 = 0
    36 = 1
    37 = 0
    38 = 2
    39 = 3
    40 = 0
    41 = 1
    42 = 2
    44 = 1
    45 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
